$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A527").EntireRow.Insert()

$ws.Cells.Item(527, 1).Value = 9
$ws.Cells.Item(527, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(527, 3).Value = "Metropolitana"
$ws.Cells.Item(527, 4).Value = 45209
$ws.Cells.Item(527, 5).Value = 13
$ws.Cells.Item(527, 6).Value = 100112039
$ws.Cells.Item(527, 7).Value = "Ciboulette"
$ws.Cells.Item(527, 8).Value = "Sin especificar"
$ws.Cells.Item(527, 9).Value = "Primera"
$ws.Cells.Item(527, 10).Value = 340
$ws.Cells.Item(527, 11).Value = 1000
$ws.Cells.Item(527, 12).Value = 1200
$ws.Cells.Item(527, 13).Value = 1100
$ws.Cells.Item(527, 14).Value = "`$/docena de atados"
$ws.Cells.Item(527, 15).Value = "Región Metropolitana"
$ws.Cells.Item(527, 16).Value = 367
$ws.Cells.Item(527, 17).Value = 3
$ws.Cells.Item(527, 18).Value = "Hortaliza"

$v = $ws.Cells.Item(527, 4).Value()
Write-Host "D527:" $v
$v2 = $ws.Cells.Item(528, 4).Value()
Write-Host "D528:" $v2
$dim = $ws.UsedRange.Rows.Count
Write-Host "Rows:" $dim
